# Auto-generated edit script: update TPM-derived NATMI metrics
# for Spp1-Itga9.xlsx (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7476426666666667
$ws.Range("H2").Value = 2.242928
$ws.Range("I2").Value = 0.001581772089386036
$ws.Range("J2").Value = 0.001581772089386036
$ws.Range("M2").Value = 4.717738333333333
$ws.Range("N2").Value = 14.153215
$ws.Range("O2").Value = 0.2002263444295212
$ws.Range("P2").Value = 0.2002263444295212
$ws.Range("Q2").Value = 3.527182468168889
$ws.Range("R2").Value = 31.74464221352
$ws.Range("S2").Value = 0.0003167124431784117
$ws.Range("T2").Value = 0.0003167124431784119
$ws.Range("G3").Value = 0.7476426666666667
$ws.Range("H3").Value = 2.242928
$ws.Range("I3").Value = 0.001581772089386036
$ws.Range("J3").Value = 0.001581772089386036
$ws.Range("N3").Value = 9.228847
$ws.Range("O3").Value = 0.1305610278731266
$ws.Range("P3").Value = 0.1305610278731266
$ws.Range("Q3").Value = 2.299959927112889
$ws.Range("R3").Value = 20.699639344016
$ws.Range("S3").Value = 0.0002065177898512639
$ws.Range("T3").Value = 0.000206517789851264
$ws.Range("G4").Value = 0.7476426666666667
$ws.Range("H4").Value = 2.242928
$ws.Range("I4").Value = 0.001581772089386036
$ws.Range("J4").Value = 0.001581772089386036
$ws.Range("M4").Value = 0.6908423333333333
$ws.Range("N4").Value = 2.072527
$ws.Range("O4").Value = 0.02932015834857891
$ws.Range("P4").Value = 0.02932015834857891
$ws.Range("Q4").Value = 0.5165032043395555
$ws.Range("R4").Value = 4.648528839056
$ws.Range("S4").Value = 0.00004637780813216108
$ws.Range("T4").Value = 0.00004637780813216109
$ws.Range("G5").Value = 0.7476426666666667
$ws.Range("H5").Value = 2.242928
$ws.Range("I5").Value = 0.001581772089386036
$ws.Range("J5").Value = 0.001581772089386036
$ws.Range("M5").Value = 15.077163
$ws.Range("N5").Value = 45.231489
$ws.Range("O5").Value = 0.6398924693487733
$ws.Range("P5").Value = 0.6398924693487733
$ws.Range("Q5").Value = 11.272330351088
$ws.Range("R5").Value = 101.450973159792
$ws.Range("S5").Value = 0.001012164048224199
$ws.Range("T5").Value = 0.001012164048224199
$ws.Range("I6").Value = 0.02590993131491687
$ws.Range("J6").Value = 0.02590993131491688
$ws.Range("M6").Value = 4.717738333333333
$ws.Range("N6").Value = 14.153215
$ws.Range("O6").Value = 0.2002263444295212
$ws.Range("P6").Value = 0.2002263444295212
$ws.Range("Q6").Value = 57.77637378903778
$ws.Range("R6").Value = 519.9873641013401
$ws.Range("S6").Value = 0.005187850831605782
$ws.Range("T6").Value = 0.005187850831605784
$ws.Range("I7").Value = 0.02590993131491687
$ws.Range("J7").Value = 0.02590993131491688
$ws.Range("N7").Value = 9.228847
$ws.Range("O7").Value = 0.1305610278731266
$ws.Range("P7").Value = 0.1305610278731266
$ws.Range("Q7").Value = 37.67407715588578
$ws.Range("R7").Value = 339.066694402972
$ws.Range("S7").Value = 0.003382827264597658
$ws.Range("T7").Value = 0.003382827264597658
$ws.Range("I8").Value = 0.02590993131491687
$ws.Range("J8").Value = 0.02590993131491688
$ws.Range("M8").Value = 0.6908423333333333
$ws.Range("N8").Value = 2.072527
$ws.Range("O8").Value = 0.02932015834857891
$ws.Range("P8").Value = 0.02932015834857891
$ws.Range("Q8").Value = 8.460487220739111
$ws.Range("R8").Value = 76.144384986652
$ws.Range("S8").Value = 0.0007596832889541662
$ws.Range("T8").Value = 0.0007596832889541663
$ws.Range("I9").Value = 0.02590993131491687
$ws.Range("J9").Value = 0.02590993131491688
$ws.Range("M9").Value = 15.077163
$ws.Range("N9").Value = 45.231489
$ws.Range("O9").Value = 0.6398924693487733
$ws.Range("P9").Value = 0.6398924693487733
$ws.Range("Q9").Value = 184.644366350596
$ws.Range("R9").Value = 1661.799297155364
$ws.Range("S9").Value = 0.01657956992975927
$ws.Range("T9").Value = 0.01657956992975927
$ws.Range("G10").Value = 18.93023433333333
$ws.Range("H10").Value = 56.79070299999999
$ws.Range("I10").Value = 0.04005030430848061
$ws.Range("J10").Value = 0.04005030430848062
$ws.Range("M10").Value = 4.717738333333333
$ws.Range("N10").Value = 14.153215
$ws.Range("O10").Value = 0.2002263444295212
$ws.Range("P10").Value = 0.2002263444295212
$ws.Range("Q10").Value = 89.30789217334943
$ws.Range("R10").Value = 803.7710295601448
$ws.Range("S10").Value = 0.008019126024976975
$ws.Range("T10").Value = 0.008019126024976977
$ws.Range("G11").Value = 18.93023433333333
$ws.Range("H11").Value = 56.79070299999999
$ws.Range("I11").Value = 0.04005030430848061
$ws.Range("J11").Value = 0.04005030430848062
$ws.Range("N11").Value = 9.228847
$ws.Range("O11").Value = 0.1305610278731266
$ws.Range("P11").Value = 0.1305610278731266
$ws.Range("Q11").Value = 58.23474544549344
$ws.Range("R11").Value = 524.1127090094409
$ws.Range("S11").Value = 0.005229008897146739
$ws.Range("T11").Value = 0.00522900889714674
$ws.Range("G12").Value = 18.93023433333333
$ws.Range("H12").Value = 56.79070299999999
$ws.Range("I12").Value = 0.04005030430848061
$ws.Range("J12").Value = 0.04005030430848062
$ws.Range("M12").Value = 0.6908423333333333
$ws.Range("N12").Value = 2.072527
$ws.Range("O12").Value = 0.02932015834857891
$ws.Range("P12").Value = 0.02932015834857891
$ws.Range("Q12").Value = 13.07780725738678
$ws.Range("R12").Value = 117.700265316481
$ws.Range("S12").Value = 0.001174281264233424
$ws.Range("T12").Value = 0.001174281264233424
$ws.Range("G13").Value = 18.93023433333333
$ws.Range("H13").Value = 56.79070299999999
$ws.Range("I13").Value = 0.04005030430848061
$ws.Range("J13").Value = 0.04005030430848062
$ws.Range("M13").Value = 15.077163
$ws.Range("N13").Value = 45.231489
$ws.Range("O13").Value = 0.6398924693487733
$ws.Range("P13").Value = 0.6398924693487733
$ws.Range("Q13").Value = 285.4142286718629
$ws.Range("R13").Value = 2568.728058046766
$ws.Range("S13").Value = 0.02562788812212347
$ws.Range("T13").Value = 0.02562788812212348
$ws.Range("G14").Value = 440.7369333333333
$ws.Range("H14").Value = 1322.2108
$ws.Range("I14").Value = 0.9324579922872165
$ws.Range("J14").Value = 0.9324579922872166
$ws.Range("M14").Value = 4.717738333333333
$ws.Range("N14").Value = 14.153215
$ws.Range("O14").Value = 0.2002263444295212
$ws.Range("P14").Value = 0.2002263444295212
$ws.Range("Q14").Value = 2079.281525302444
$ws.Range("R14").Value = 18713.533727722
$ws.Range("S14").Value = 0.18670265512976
$ws.Range("T14").Value = 0.1867026551297601
$ws.Range("G15").Value = 440.7369333333333
$ws.Range("H15").Value = 1322.2108
$ws.Range("I15").Value = 0.9324579922872165
$ws.Range("J15").Value = 0.9324579922872166
$ws.Range("N15").Value = 9.228847
$ws.Range("O15").Value = 0.1305610278731266
$ws.Range("P15").Value = 0.1305610278731266
$ws.Range("Q15").Value = 1355.831241660844
$ws.Range("R15").Value = 12202.4811749476
$ws.Range("S15").Value = 0.1217426739215309
$ws.Range("T15").Value = 0.1217426739215309
$ws.Range("G16").Value = 440.7369333333333
$ws.Range("H16").Value = 1322.2108
$ws.Range("I16").Value = 0.9324579922872165
$ws.Range("J16").Value = 0.9324579922872166
$ws.Range("M16").Value = 0.6908423333333333
$ws.Range("N16").Value = 2.072527
$ws.Range("O16").Value = 0.02932015834857891
$ws.Range("P16").Value = 0.02932015834857891
$ws.Range("Q16").Value = 304.4797314101777
$ws.Range("R16").Value = 2740.317582691599
$ws.Range("S16").Value = 0.02733981598725916
$ws.Range("T16").Value = 0.02733981598725916
$ws.Range("G17").Value = 440.7369333333333
$ws.Range("H17").Value = 1322.2108
$ws.Range("I17").Value = 0.9324579922872165
$ws.Range("J17").Value = 0.9324579922872166
$ws.Range("M17").Value = 15.077163
$ws.Range("N17").Value = 45.231489
$ws.Range("O17").Value = 0.6398924693487733
$ws.Range("P17").Value = 0.6398924693487733
$ws.Range("Q17").Value = 6645.062583986799
$ws.Range("R17").Value = 59805.56325588118
$ws.Range("S17").Value = 0.5966728472486664
$ws.Range("T17").Value = 0.5966728472486664
